$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observations")

# New row of data (row 3)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 45141
$ws.Range("B3").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("C3").Value = "Send / Request "
$ws.Range("D3").Value = "Add Card to Reload"
$ws.Range("F3").Value = "Debit,Credit card address is  not auto populating without existing card but addressed already saved"
$ws.Range("E3").Value = "Address fields"

# Column widths (stored OOXML width = ColumnWidth + 5/6; compensate so
# the saved file shows width="15.5" / width="61.5" like the target)
$ws.Columns.Item(3).ColumnWidth = 44.0 / 3.0
$ws.Columns.Item(4).ColumnWidth = 44.0 / 3.0
$ws.Columns.Item(6).ColumnWidth = 182.0 / 3.0

# Selection
$ws.Range("D3").Select()
